$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("EMPLOYEE DTR")

# Clear the "01:00:00" values from the OFFICIAL BUSINESS DEPARTURE (K)
# and OFFICIAL BUSINESS ARRIVAL (N) columns for rows 5-15 and 17-18.
$rows = @(5,6,7,8,9,10,11,12,13,14,15,17,18)
foreach ($r in $rows) {
    $ws.Cells.Item($r, 11).ClearContents()
    $ws.Cells.Item($r, 14).ClearContents()
}
